$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F5 becomes a text value "NA" instead of the date 45534
$ws.Range("F5").Value = "NA"

# F6 date changes from 45534 to 45595 (2024-10-30)
$ws.Range("F6").Value = [DateTime]::FromOADate(45595)

# Update the selection on the sheet to F6 (single cell)
$ws.Range("F6").Select()
